$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new order-item row (row 56) after the existing last row (row 55)
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 19
$ws.Range("C56").Value = "CAND234"
$ws.Range("D56").Value = 2
$ws.Range("E56").Value = 1793.72197309417
$ws.Range("F56").Value = 0
